# Weekly update: insert a new price record as row 49, pushing the
# existing rows 49-69 down to 50-70 (the sheet keeps a rolling history,
# newest entry on top of the dated list).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 49; this shifts old rows 49:69
# down to 50:70 and keeps their formatting/values intact.
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row 49 with this week's data.
$ws.Cells.Item(49, 1).Value = 8
$ws.Cells.Item(49, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(49, 3).Value = "Coquimbo"
$ws.Cells.Item(49, 4).Value = 44846
$ws.Cells.Item(49, 5).Value = 4
$ws.Cells.Item(49, 6).Value = 100114007
$ws.Cells.Item(49, 7).Value = "Jengibre"
$ws.Cells.Item(49, 8).Value = "Sin especificar"
$ws.Cells.Item(49, 9).Value = "Primera"
$ws.Cells.Item(49, 10).Value = 400
$ws.Cells.Item(49, 11).Value = 13500
$ws.Cells.Item(49, 12).Value = 14000
$ws.Cells.Item(49, 13).Value = 13750
$ws.Cells.Item(49, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(49, 15).Value = "Perú"
$ws.Cells.Item(49, 16).Value = 1058
$ws.Cells.Item(49, 17).Value = 13
$ws.Cells.Item(49, 18).Value = "Hortaliza"
